$p = $ppt.ActivePresentation
$s = $p.Slides.Item(21)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Insert the 6 new level-1 sub-bullet paragraphs right after paragraph 1 ---
# (done first, using the final text verbatim, so no later edit causes the
#  diff-style run engine to fragment these single-run paragraphs)
$para1 = $tr.Paragraphs(1, 1)
$newText = "`rgit checkout repo`rgit checkout -b branchname`rAdd/Edit/Delete files`rgit add –A .`rgit commit -m “change comment”`rgit push"
$para1.InsertAfter($newText)

# --- Bump the indent level of the 6 newly inserted paragraphs (now
#     paragraphs 2-7) to level 2 (=> lvl="1" in the OOXML). ---
$tr = $sh.TextFrame.TextRange
for ($i = 2; $i -le 7; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.IndentLevel = 2
}

# --- Split the "git checkout -b branchname" paragraph (paragraph 3) into
#     two runs: "git checkout -b " and "branchname" ---
$tr = $sh.TextFrame.TextRange
$para3 = $tr.Paragraphs(3, 1)
$branchPart = $para3.Characters(17, 10)
$branchPart.Text = "branchname"

# --- Finally, split paragraph 1 ("Basic GIT commands") into two runs:
#     "Basic GIT " and "workflow and commands" ---
$tr = $sh.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)
$secondPart = $para1.Characters(11, 8)
$secondPart.Text = "workflow and commands"
